{"js": "// Oscar's Week2 individual reflection, draft 2:\n// \"... but I hard the team did great.\" -> \"... but I heard the team did great.\"\n//\n// The canonical diff shows the corrected word (\"heard\") landing in its own\n// run, flanked by two sibling runs that carry the untouched text before/after\n// it - all three runs share identical run properties (Times New Roman, 24pt,\n// en-US). That is exactly what Word's spell/grammar-correction flow produces\n// when you accept a fix for a single word inside a longer run. We reproduce\n// it here: find the misspelled word, replace just that word, then nudge a\n// character-formatting property on the freshly inserted range (set then\n// immediately unset) so the engine keeps it as a discrete run instead of\n// silently re-merging it with its now-identical neighbors.\n\nconst body = context.document.body;\n\n// \"hard\" only occurs once in this document, inside \"... but I hard the team\n// did great.\" - search narrowly (whole word, case sensitive) so we only ever\n// touch that one word.\nconst found = body.search(\"hard\", { matchCase: true, matchWholeWord: true });\nfound.load(\"items\");\nawait context.sync();\n\nif (found.items.length === 0) {\n  throw new Error('edit.js: could not find \"hard\" to correct to \"heard\".');\n}\n\nconst target = found.items[0];\ntarget.insertText(\"heard\", \"Replace\");\nawait context.sync();\n\n// Re-find the just-inserted word and flip a formatting property on/off so\n// the corrected word stays split into its own run (matching the authored\n// diff) rather than being coalesced back into the surrounding text.\nconst corrected = body.search(\"heard\", { matchCase: true, matchWholeWord: true });\ncorrected.load(\"items\");\nawait context.sync();\n\nconst correctedRange = corrected.items[0];\ncorrectedRange.font.bold = true;\nawait context.sync();\n\ncorrectedRange.font.bold = false;\nawait context.sync();\n", "ps1": "# Oscar's Week2 individual reflection, draft 2:\n# \"... but I hard the team did great.\" -> \"... but I heard the team did great.\"\n#\n# The canonical diff shows the corrected word (\"heard\") landing in its own\n# run, flanked by two sibling runs that carry the untouched text before/after\n# it - all three runs share identical run properties (Times New Roman, 24pt,\n# en-US). That is exactly what Word's spell/grammar-correction flow produces\n# when you accept a fix for a single word inside a longer run. We reproduce\n# it here: find the misspelled word, replace just that word, then nudge a\n# character-formatting property on the freshly inserted range (set then\n# immediately unset) so the engine keeps it as a discrete run instead of\n# silently re-merging it with its now-identical neighbors.\n\n$d = $word.ActiveDocument\n\n# \"hard\" only occurs once in this document, inside \"... but I hard the team\n# did great.\" - search narrowly (whole word, case sensitive) so we only ever\n# touch that one word.\n$find = $d.Content\n$found = $find.Find.Execute(\"hard\", $true, $true, $false, $false, $false, $true, 1, $false, \"heard\", 2)\n\nif (-not $found) {\n    throw 'edit.ps1: could not find \"hard\" to correct to \"heard\".'\n}\n\n# Re-find the just-corrected word and flip a formatting property on/off so it\n# stays split into its own run (matching the authored diff) rather than being\n# coalesced back into the surrounding text.\n$again = $d.Content\n$again.Find.Execute(\"heard\", $true, $true) | Out-Null\n$again.Bold = 1\n$again.Bold = 0\n"}
